# Commit message: "repull data, push all data, mean calculation"
# The dSF column (column F) values were repulled/recalculated for most rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 2
    6  = -2
    7  = -5
    8  = -1
    9  = -1
    10 = 5
    11 = 3
    12 = -3
    13 = -5
    14 = 2
    16 = -2
    17 = -1
    18 = 5
    19 = -2
    20 = -2
    21 = 2
    22 = 6
    23 = -2
    24 = 3
    25 = -4
    26 = -2
    27 = -4
    28 = -2
    30 = -1
    31 = -6
    32 = 2
    34 = -1
    36 = -9
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
